$d = $word.ActiveDocument
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:t>ANDREW MERCER</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">One throne. Hundreds of thousands of heirs. Some of them live lavishly with a near infinite wealth, claiming their blood is purer and that it raises their place in the line of succession. Others don’t care and simply want to live their lives. </w:t>
      </w:r>
      <w:r>
        <w:t>Only a f</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ew think about the people and the clans, few want to rule for reasons other than personal glory and power. One of these few is Andrew Mercer. He sees the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerite</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Clans have begun to decline in power and prestige, as all great empires and kingdoms eventually do. Most have given up and seek to ride out the last days in what’s left of their wealth and power. Andrew Mercer seeks to restore the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerite</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Clans to their former glory and possibly to a level beyond that. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">Like the great warlord of old, Andrew’s life has been steeped in struggle. The original Mercer relished struggle because of the great triumph that comes afterwards, but Andrew did not have such joy with his trials. He was born with a poor link to the throne. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>Almost all his close friends and family were extinguished by one of many “Purges” by more high ranking heirs to the throne.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">He barely survived during his service in the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinari</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Survival War. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">And the majority of </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">citizens of the Spacer Alliance that are </w:t>
      </w:r>
      <w:r>
        <w:t>non-</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerite</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> are prejudiced against him, for they fear and detest his infamous last name. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t>So to display his ability to rule he has taken on mission aft</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">er mission to serve </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerite</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Clans. Over the years many </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">minor </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">clans have grown fond of him as he helped them over and over again. And in time, he gathered a faithful </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>warband</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> around him, some from Great Clans, others from minor clans, and even some from Clan Mercer itself. Over time he has gained</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> some</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> recognition for his deeds,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> but it is not enough if he wants to take the throne. He has set his eyes on a new target now; a massive criminal network called “the Family” is disrupting trade in the Spacer Alliance. Putting a stop to their crime will be a great victory for the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerites</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Vinari</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve"> Survival War</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">Andrew Mercer spent some of his years fighting in the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinari</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Survival War. He was sent by higher ranking Mercers because they believed he would die, as he had little to no real combat </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>experience</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> before then. Many lower ranking Mercers were sent to fight in the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinari</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Survival War because the odds of surviving a term of service on the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinnas</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Front was incredibly low. Andrew not only survived through his term, he also earned himself a </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>warband</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> as he was promoted repeatedly. When he left the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinnas</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Front, he was able to take most of his </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>warband</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> with him.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">During his time on the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinnas</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Front, Andrew realized the perilous state of the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerite</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Clans. He saw the tremendous effort of the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinari</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> peoples in their attempt to push back against their enemy. He befriended a number of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinari</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, both </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Angelfolk</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> and Wingless, and discovered that their united focus was something that the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerites</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> were lacking. With the fall of the Densen Clan, the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerites</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> were no longer an economic powerhouse as their main trading partner was gone. Additionally, due to poor planning but good intentions, participating in the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinari</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Survival War was stretching </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerite</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>warbands</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> too thinly across space. The morale during the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinari</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Survival War was low, and many </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerites</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> who weren’t on the front were worried about </w:t>
      </w:r>
      <w:r>
        <w:t>the implications of losing that war</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. There was even talk of division as some </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">minor </w:t>
      </w:r>
      <w:r>
        <w:t>clans expressed a desire to leave o</w:t>
      </w:r>
      <w:r>
        <w:t>r rebel against their Great Clan</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve">The </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Warband</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:tab/>
        <w:t xml:space="preserve">Andrew’s </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>warband</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> consist of just about one hundred Spacers. They come from a variety of Great Clans and minor clans. A few are also from non-</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerite</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> clans. Andrew proves to be a tactful leader in his ability to rally so many Spacers of varying opinions and cultures. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">He manages to lead them with few words as they fear and adore him and do not need too much inspiration to enact his orders; they have seen already what he is capable of. For more on important members of his </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>warband</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, see the related document in the same folder as this document.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Personality</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Andrew has a certain focus that is incomparable to most people. When he wants to get something done, he is willing to do anything in his power to obtain it except sacrificing his </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>warband</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>. In this sense, h</w:t>
      </w:r>
      <w:r>
        <w:t>is</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>warband</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> is both his greatest strength and his greatest weakness. Through a subtle charisma and his willingness to acknowledge the everyman, Andrew has managed to captivate his </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>warband</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">. However, they are tied so closely that their pain is his, and his pain is theirs. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">Just like his namesake, Andrew also walks between the line of civility and savagery. Depending on what is more likely to get the results he wants, he may choose to talk or he may choose violence. He is especially adept at the latter due to his time on the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinnas</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Front. However, when he does choose the peaceful path, he often attempts to empathize with the other</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> party</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Lastly, Andrew has a sense of discipline. He keeps himself combat ready by training or drilling on a daily basis. He attempts to match this discipline with honor but struggles to do so. Rather, he believes in the concept of survival over anything else, as he saw its importance and power firsthand on the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vinnas</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Front. Additionally, honor has less value to him as he witnessed the dishonorable actions of many high-ranking </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerites</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> who believe themselves to be extremely honorable. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Summary</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">The few that have met Andrew Mercer agree that he should one day rule the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Mercerites</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">. In a time of chaos and confusion, he would make an excellent leader. However, this does not matter to the protagonist, who is seen as an enemy by Andrew. The dealings that will occur between Andrew and the protagonist will have an impact on Andrew, the protagonist, and the overall story of the game as well. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p></w:body></pkg:xmlData></pkg:part></pkg:package>
'@
$r = $d.Content
$r.InsertXML($xml)
